# The presentation's slide master currently uses the "Integral" theme
# (ppt/theme/theme1.xml). The target edit swaps the theme applied to the
# deck for the default "Office Theme" color palette (the theme that used
# to sit, unused, behind the Notes Master). We reproduce the
# visually-significant part of that swap -- the 12 theme colors -- via
# the Theme Color Scheme that PowerPoint exposes on every slide (it is
# backed by the one Slide Master / Theme that drives the whole deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB integers (0xBBGGRR), matching the target OOXML
# srgbClr values: 000000, FFFFFF, 44546A, E7E6E6, 5B9BD5, ED7D31,
# A5A5A5, FFC000, 4472C4, 70AD47, 0563C1, 954F72.
$officeColors = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
